$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column S (19th column) - this shifts all
# subsequent columns (data, styles, column-width specs, dimension) one to
# the right automatically, matching the "Sub brand" column insertion.
$ws.Columns("S:S").Insert()

# New header for the inserted column.
$ws.Range("S1").Value = "Sub brand"

# Re-point the autofilter so it covers the newly inserted column
# (old range A1:AO54 -> new range A1:AP54).
$ws.AutoFilterMode = $false
$ws.Range("A1:AP54").AutoFilter()

# Keep the two _FilterDatabase defined names in sync with the wider range.
foreach ($n in $wb.Names) {
    $n.RefersTo = "=Cinema!`$A`$1:`$AP`$54"
}

# Turn off iterative calculation (author cleared calcPr's iterateDelta).
$excel.IterativeCalculation = $false

# Move the active selection to the newly inserted column's second row.
$ws.Range("S2").Select()
